$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$cellValues = @{
    "M9" = 32
    "H9" = 1616.3334
    "N9" = -4913
    "L9" = 4575
    "K9" = 137
    "J9" = 4575
    "I9" = 137
    "M17" = -3331.8
    "K17" = 3499.8
    "I17" = 1166.6
    "H17" = 1479.7059
    "K43" = 6250
    "L43" = 13486.286
    "J43" = 13486.286
    "M43" = -6181
    "N43" = -13624.286
    "I43" = 6250
    "K58" = 5537.4
    "J58" = 7754.25
    "I58" = 1845.8
    "N58" = -23562.75
    "L58" = 23262.75
    "M58" = -5387.4
    "H58" = 5481.769
    "N112" = -8131.549999999999
    "J112" = 1971.85
    "L112" = 5915.549999999999
    "H112" = 1939.8572
    "L138" = 15420.522
    "J138" = 5140.174
    "N138" = -25700.522
    "H138" = 5142.041
}
foreach ($key in $cellValues.Keys) {
    $ws.Range($key).Value = $cellValues[$key]
}

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$cellValues = @{
    "J32" = 46338
    "N32" = -46912
    "K32" = 7769.885
    "L32" = 46338
    "M32" = -7482.885
    "H32" = 9873.6
    "I32" = 7769.885
    "I61" = 4607.524
    "N61" = -8559.333500000001
    "H61" = 5391.4814
    "J61" = 8135.3335
    "M61" = -4395.524
    "K61" = 4607.524
    "L61" = 8135.3335
    "I63" = 3418.3333
    "H63" = 7437.75
    "J63" = 8777.556
    "L63" = 8777.556
    "N63" = -10149.556
    "K63" = 3418.3333
    "M63" = -2732.3333
    "N66" = -50751.78
    "J66" = 8777.556
    "H66" = 7437.75
    "K66" = 17091.6665
    "L66" = 43887.78
    "I66" = 3418.3333
    "M66" = -13659.6665
    "I74" = 2582.2
    "M74" = -1708.2
    "L74" = 5794
    "H74" = 4334.091
    "J74" = 5794
    "K74" = 2582.2
    "N74" = -7542
    "J77" = 5794
    "H77" = 4334.091
    "M77" = -8543
    "I77" = 2582.2
    "L77" = 28970
    "K77" = 12911
    "N77" = -37706
    "L97" = 1397.6
    "M97" = -378.3333
    "N97" = -2389.6
    "I97" = 874.3333
    "H97" = 1112.1818
    "K97" = 874.3333
    "J97" = 1397.6
    "J122" = 4430.091
    "M122" = -4720.24
    "K122" = 7170.24
    "N122" = -18190.273
    "L122" = 13290.273
    "I122" = 2390.08
    "H122" = 3013.4167
    "M132" = -6126.4547
    "H132" = 3353.3157
    "N132" = -24383
    "L132" = 19323
    "J132" = 6441
    "K132" = 8656.4547
    "I132" = 2885.4849
    "L136" = 24406.0005
    "K136" = 13822.572
    "H136" = 5391.4814
    "M136" = -11272.572
    "I136" = 4607.524
    "N136" = -29506.0005
    "J136" = 8135.3335
}
foreach ($key in $cellValues.Keys) {
    $ws.Range($key).Value = $cellValues[$key]
}

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$cellValues = @{
    "J54" = 1100
    "H54" = 1416.6666
    "N54" = -2068
    "L54" = 1100
    "I59" = 25000
    "K59" = 25000
    "L59" = 0
    "J59" = 0
    "H59" = 25000
    "M86" = -2514.7368
    "H86" = 4524
    "K86" = 3637.7368
    "I86" = 3637.7368
    "I89" = 3637.7368
    "M89" = -12572.684
    "H89" = 4524
    "K89" = 18188.684
    "K94" = 2709.0557
    "L94" = 3166.6667
    "M94" = -2258.0557
    "J94" = 3166.6667
    "H94" = 2774.4285
    "I94" = 2709.0557
    "N94" = -4068.6667
    "I134" = 4185.0527
    "N134" = -20287.7139
    "L134" = 15217.7139
    "H134" = 4424
    "K134" = 12555.1581
    "M134" = -10020.1581
    "J134" = 5072.5713
}
foreach ($key in $cellValues.Keys) {
    $ws.Range($key).Value = $cellValues[$key]
}
$ws.Range("N59").ClearContents()
$ws.Range("M59").Value = -24153

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$cellValues = @{
    "L31" = 120778.445
    "N31" = -121368.445
    "I31" = 14711
    "J31" = 120778.445
    "M31" = -14416
    "K31" = 14711
    "H31" = 64953.473
    "N34" = -121182.445
    "J34" = 120778.445
    "M34" = -14509
    "K34" = 14711
    "H34" = 64953.473
    "I34" = 14711
    "L34" = 120778.445
    "J122" = 10012.5
    "M122" = -7576
    "K122" = 10026
    "N122" = -34937.5
    "L122" = 30037.5
    "I122" = 3342
    "H122" = 6899.6
}
foreach ($key in $cellValues.Keys) {
    $ws.Range($key).Value = $cellValues[$key]
}

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$cellValues = @{
    "L23" = 2199.75
    "K23" = 1279.00002
    "M23" = -1044.00002
    "I23" = 426.33334
    "H23" = 649.5454999999999
    "J23" = 733.25
    "N23" = -2669.75
    "H93" = 7496.6665
    "J93" = 4997.5
    "N93" = -18736.5
    "L93" = 14992.5
    "J122" = 3410.0454
    "M122" = -8029.6
    "K122" = 10479.6
    "N122" = -35590.4086
    "L122" = 30690.4086
    "I122" = 1164.4
    "H122" = 2708.2812
    "M124" = -15590.0005
    "I124" = 6833.3335
    "H124" = 7587.222
    "K124" = 20500.0005
    "H137" = 75187.5
    "L137" = 285238.89
    "J137" = 95079.63
    "N137" = -295438.89
    "J139" = 4579.0835
    "H139" = 3518.9656
    "L139" = 13737.2505
    "K139" = 8311.940999999999
    "I139" = 2770.647
    "N139" = -24017.2505
    "M139" = -3171.940999999999
}
foreach ($key in $cellValues.Keys) {
    $ws.Range($key).Value = $cellValues[$key]
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$cellValues = @{
    "K54" = 10000
    "I54" = 10000
    "J54" = 10666.667
    "H54" = 10500
    "N54" = -11446.667
    "L54" = 10666.667
    "L70" = 11430
    "N70" = -11970
    "H70" = 9922.096
    "J70" = 11430
    "J73" = 11430
    "L73" = 11430
    "N73" = -13302
    "H73" = 9922.096
    "N126" = -19321.25
    "I126" = 4986.45
    "J126" = 4793.75
    "H126" = 4900.8057
    "M126" = -12489.35
    "L126" = 14381.25
    "K126" = 14959.35
    "M132" = -6745.625
    "H132" = 4614.8184
    "N132" = -31088
    "L132" = 26028
    "J132" = 8676
    "K132" = 9275.625
    "I132" = 3091.875
}
foreach ($key in $cellValues.Keys) {
    $ws.Range($key).Value = $cellValues[$key]
}
$ws.Range("M54").Value = -9610

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$cellValues = @{
    "N22" = -14589.571
    "J22" = 13999.571
    "M22" = -250000765
    "H22" = 90918380
    "L22" = 13999.571
    "K22" = 250001060
    "I22" = 250001060
    "K27" = 250001060
    "H27" = 90918380
    "L27" = 13999.571
    "I27" = 250001060
    "M27" = -250000953
    "N27" = -14213.571
    "J27" = 13999.571
    "J44" = 100000
    "N44" = -100912
    "H44" = 100000
    "L44" = 100000
    "M93" = -1561.5
    "K93" = 2809.5
    "H93" = 2781.3635
    "I93" = 2809.5
    "M132" = -10961.4995
    "H132" = 5331.6113
    "N132" = -26061.5
    "L132" = 21001.5
    "J132" = 7000.5
    "K132" = 13491.4995
    "I132" = 4497.1665
}
foreach ($key in $cellValues.Keys) {
    $ws.Range($key).Value = $cellValues[$key]
}

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$cellValues = @{
    "M132" = -3044.5625
    "H132" = 3500.5
    "K132" = 5574.5625
    "I132" = 1858.1875
    "L136" = 36869.331
    "K136" = 7014.8181
    "H136" = 4470.7383
    "M136" = -4464.8181
    "I136" = 2338.2727
    "N136" = -41969.331
    "J136" = 12289.777
}
foreach ($key in $cellValues.Keys) {
    $ws.Range($key).Value = $cellValues[$key]
}
